$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels: "Conditioned" -> "Unconditioned"
$ws.Range("B1").Value = "1:4 Unconditioned"
$ws.Range("C1").Value = "4:1 Unconditioned"

# Widen column C slightly
$ws.Columns("C").ColumnWidth = 16.1640625

# Move the active selection to D7
$ws.Range("D7").Select()
